$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 10, 11, 12, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-22 18:22:26"
    $zhcn.Range("H$r").Value = "2016-08-22 18:22:20"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
